# Datapath_Control_Sheet.xlsx edit script
# Summary of the change (per commit "lab2-1: declare global control constant"):
#   - Rename the "BrUnSel" control signal to "BrSel" throughout both sheets
#     (header cells and every literal value like BrUnSel_Non/_Beq/_Bne/_Blt/_Bge).
#   - Correct the BrUn column entries for the beq/bne/blt/bge example rows
#     (J25:J28) from the stray "BrUn_En(1'b1)" back to "BrUn_Non".
#   - Fix two stray leftover "BrUn"/"BrUnSel" labels on the datapath sheet
#     (N35 -> "BrSel", N37 -> "-").
#   - Switch the active tab to "控制信号取值表" (sheet 2) and restore each
#     sheet's last-used selection.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # 数据通路表(含控制信号)
$ws2 = $wb.Worksheets.Item(2)   # 控制信号取值表

# --- Sheet 2 ("控制信号取值表"): BrUnSel -> BrSel header + value strings ---

# Column header (row 3)
$ws2.Range("K3").Value = "BrSel"

# "Non" (default / don't-care) value, used by most instruction rows
$nonRows = @(4)
foreach ($r in $nonRows) {
    $ws2.Range("K$r").Value = "BrSel_Non(3'b0)"
}
$nonRows2 = @(5,6,7,8,9,10,11,13,14,15,16,17,18,19,20,21,23,30,32)
foreach ($r in $nonRows2) {
    $ws2.Range("K$r").Value = "BrSel_Non"
}

# Branch-specific selector values (beq/bne/blt/bge example rows)
$ws2.Range("K25").Value = "BrSel_Beq(3'b1)"
$ws2.Range("K26").Value = "BrSel_Bne(3'b10)"
$ws2.Range("K27").Value = "BrSel_Blt(3'b11)"
$ws2.Range("K28").Value = "BrSel_Bge(3'b100)"

# Fix the BrUn column for the same branch rows: was mistakenly "BrUn_En(1'b1)",
# should read "BrUn_Non" like every other instruction row.
$brUnRows = @(25,26,27,28)
foreach ($r in $brUnRows) {
    $ws2.Range("J$r").Value = "BrUn_Non"
}

# --- Sheet 1 ("数据通路表(含控制信号)"): fix stray BrUn/BrUnSel labels ---
$ws1.Range("N35").Value = "BrSel"
$ws1.Range("N37").Value = "-"

# --- View state: make "控制信号取值表" the active tab, restore selections ---
$ws1.Range("N38").Select()
$ws2.Range("L16").Select()
$ws2.Activate()
